$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> @(WIID, Status, ExcelDateSerial)
$updates = @{
    2  = @(809944, "Open", 43232)
    3  = @(619534, "Open", 43143)
    4  = @(532044, "Open", 42903)
    5  = @(612294, "Open", 42895)
    6  = @(640094, "Open", 42996)
    7  = @(655994, "Open", 43171)
    8  = @(640314, "Open", 43294)
    9  = @(637964, "Open", 43227)
    10 = @(193134, "Open", 43193)
    11 = @(798064, "Open", 42913)
    12 = @(442854, "Open", 42995)
    13 = @(563864, "Open", 43483)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $wiid = $vals[0]
    $status = $vals[1]
    $dateSerial = $vals[2]

    $ws.Cells.Item($row, 2).Value = $wiid
    $ws.Cells.Item($row, 5).Value = $status
    $ws.Cells.Item($row, 6).Value = $dateSerial
}
